$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1,1).Range.Text = "59 x 11" + $nl + "  1    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
$t.Cell(1,2).Range.Text = "78 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "8|    |"
$t.Cell(1,3).Range.Text = "71 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "7|    |" + $nl + "1|    |"
$t.Cell(2,1).Range.Text = "97 x 10" + $nl + "  1    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "7|    |"
$t.Cell(2,2).Range.Text = "95 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "9|    |" + $nl + "5|    |"
$t.Cell(2,3).Range.Text = "88 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "8|    |"
$t.Cell(3,1).Range.Text = "61 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "1|    |"
$t.Cell(3,2).Range.Text = "78 x 65" + $nl + "  6    5" + $nl + "  ----" + $nl + "7|    |" + $nl + "8|    |"
$t.Cell(3,3).Range.Text = "61 x 24" + $nl + "  2    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "1|    |"
$t.Cell(4,1).Range.Text = "35 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "3|    |" + $nl + "5|    |"
$t.Cell(4,2).Range.Text = "93 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$t.Cell(4,3).Range.Text = "45 x 55" + $nl + "  5    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "5|    |"
$t.Cell(5,1).Range.Text = "12 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"
$t.Cell(5,2).Range.Text = "58 x 60" + $nl + "  6    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
$t.Cell(5,3).Range.Text = "86 x 12" + $nl + "  1    2" + $nl + "  ----" + $nl + "8|    |" + $nl + "6|    |"
